# feat: add 2022-Q1 data
#
# 1) Insert a new "2022-Q1" sheet (a duplicate of "总计"'s formatting) between
#    "2021-Q3" and "总计", populated with the Q1-2022 holdings detail.
# 2) Prepend a "2022-Q1" summary row to the "总计" sheet, pushing the existing
#    "2021-Q3" row down one slot.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: new "2022-Q1" worksheet
# ---------------------------------------------------------------------------

$wsTotal = $wb.Worksheets.Item("总计")

# Duplicate "总计" (same header/column styling as the other per-quarter sheet)
# and drop the copy immediately in front of it, so the tab order becomes
# 2021-Q3, 2022-Q1, 总计.
$wsTotal.Copy($wsTotal, $null)
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q1"

# The copy still holds 总计's own B1:D2 block - wipe it before filling in
# the new sheet's data.
$ws.Range("A1:D2").ClearContents()

# Extend the header/column formatting that came along with the copy so it
# covers the wider B:H layout used by the per-quarter detail sheets.
$ws.Range("B1:D1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Columns B:G hold numeric-looking figures that must stay text (fund codes
# with leading zeros, percentages, etc.), so format them as text first.
$ws.Range("B2:G5").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "513060"
$ws.Range("C2").Value = "博时恒生医疗保健ETF（QDII）"
$ws.Range("D2").Value = "23.35"
$ws.Range("E2").Value = "98.89"
$ws.Range("F2").Value = "2.77"
$ws.Range("G2").Value = "0.6468"
$ws.Range("H2").Value = 8

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "513700"
$ws.Range("C3").Value = "鹏华中证港股通医药卫生综合交易型开放式指数证券投资基金"
$ws.Range("D3").Value = "3.24"
$ws.Range("E3").Value = "93.11"
$ws.Range("F3").Value = "2.13"
$ws.Range("G3").Value = "0.0690"
$ws.Range("H3").Value = 10

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "159892"
$ws.Range("C4").Value = "华夏恒生香港上市生物科技ETF（QDII）"
$ws.Range("D4").Value = "1.51"
$ws.Range("E4").Value = "99.03"
$ws.Range("F4").Value = "3.02"
$ws.Range("G4").Value = "0.0456"
$ws.Range("H4").Value = 9

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "004098"
$ws.Range("C5").Value = "前海开源港股通股息率50强股票"
$ws.Range("D5").Value = "0.34"
$ws.Range("E5").Value = "88.92"
$ws.Range("F5").Value = "2.21"
$ws.Range("G5").Value = "0.0075"
$ws.Range("H5").Value = 9

# ---------------------------------------------------------------------------
# Part 2: prepend the 2022-Q1 summary row on "总计"
# ---------------------------------------------------------------------------

# Re-resolve by name: the sheet collection shifted when the copy above was
# inserted, so the old $wsTotal handle (bound by the pre-copy index) would
# now point at the new "2022-Q1" sheet instead.
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Rows(2).Insert()
$wsTotal.Rows(2).ClearFormats()

# Restore column-A styling (matches the sibling cell below) on the new row.
$wsTotal.Range("A3").Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 4
$wsTotal.Range("D2").Value = 0.77

# The pushed-down "2021-Q3" row becomes index 1.
$wsTotal.Range("A3").Value = 1
